# Update cryptocurrency price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking price
# strings (e.g. "1.00", "2.44") are written as text, matching the
# original inline-string cell content instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.068.56"
$ws.Range("D3").Value = "2.305.91"
$ws.Range("D5").Value = "300.65"
$ws.Range("D6").Value = "98.47"
$ws.Range("D7").Value = "0.521"
$ws.Range("D10").Value = "35.67"
$ws.Range("D13").Value = "17.91"
$ws.Range("D15").Value = "2.663.34"
$ws.Range("D16").Value = "2.269.91"
$ws.Range("D18").Value = "42.979.91"
$ws.Range("D19").Value = "13.41"
$ws.Range("D22").Value = "68.40"
$ws.Range("D23").Value = "239.47"
$ws.Range("D25").Value = "0.999"
$ws.Range("D26").Value = "2.44"
$ws.Range("D27").Value = "24.77"
$ws.Range("D28").Value = "167.45"
$ws.Range("D31").Value = "33.38"
$ws.Range("D33").Value = "1.00"
$ws.Range("D34").Value = "4.83"
$ws.Range("D35").Value = "18.12"
$ws.Range("D42").Value = "2.007.77"
$ws.Range("D44").Value = "2.15"
$ws.Range("D45").Value = "10.06"
$ws.Range("D46").Value = "17.38"
$ws.Range("D48").Value = "54.42"
$ws.Range("D49").Value = "2.530.16"
$ws.Range("D50").Value = "73.91"

# Restore the default (Normal) style on column D so no stray number
# formatting is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"

# Volume(1h) percentage strings keep their two leading/trailing spaces.
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  +4.12%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("E19").Value = "  +7.68%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("E30").Value = "  -9.54%  "
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +4.07%  "
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +5.94%  "
$ws.Range("E51").Value = "  +1.61%  "
